$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 41666780
$ws.Range("I2").Value = 62500070
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 62500070
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = -62499957
$ws.Range("N2").Value = -426

# Row 4
$ws.Range("H4").Value = 7692672.5
$ws.Range("I4").Value = 11111428
$ws.Range("J4").Value = 473.25
$ws.Range("K4").Value = 11111428
$ws.Range("L4").Value = 473.25
$ws.Range("M4").Value = -11111314
$ws.Range("N4").Value = -701.25

# Row 19
$ws.Range("H19").Value = 430.2414
$ws.Range("I19").Value = 395.07144
$ws.Range("J19").Value = 463.06668
$ws.Range("K19").Value = 395.07144
$ws.Range("L19").Value = 463.06668
$ws.Range("M19").Value = -220.07144
$ws.Range("N19").Value = -813.06668

# Row 29
$ws.Range("H29").Value = 4200
$ws.Range("I29").Value = 3900
$ws.Range("J29").Value = 4500
$ws.Range("K29").Value = 11700
$ws.Range("L29").Value = 13500
$ws.Range("M29").Value = -11419
$ws.Range("N29").Value = -14062

# Row 31
$ws.Range("H31").Value = 1628.5714
$ws.Range("I31").Value = 280
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 840
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = -610
$ws.Range("N31").Value = -15460

# Row 32
$ws.Range("H32").Value = 600
$ws.Range("I32").Value = 600
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 600
$ws.Range("M32").Value = -274
$ws.Range("N32").Value = -1252

# Row 38
$ws.Range("H38").Value = 658.125
$ws.Range("I38").Value = 658.125
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1974.375
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -1602.375

# Row 39
$ws.Range("H39").Value = 107.2
$ws.Range("I39").Value = 112
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 336
$ws.Range("L39").Value = 300
$ws.Range("M39").Value = -40
$ws.Range("N39").Value = -892

# Row 40
$ws.Range("H40").Value = 11112906
$ws.Range("I40").Value = 1779.0883
$ws.Range("J40").Value = 45456390
$ws.Range("K40").Value = 1779.0883
$ws.Range("L40").Value = 45456390
$ws.Range("M40").Value = -1604.0883
$ws.Range("N40").Value = -45456740

# Row 41
$ws.Range("H41").Value = 669.2
$ws.Range("I41").Value = 850
$ws.Range("J41").Value = 624
$ws.Range("K41").Value = 850
$ws.Range("L41").Value = 624
$ws.Range("M41").Value = -410
$ws.Range("N41").Value = -1504

# Row 42
$ws.Range("H42").Value = 357.55554
$ws.Range("I42").Value = 61.5
$ws.Range("J42").Value = 594.4
$ws.Range("K42").Value = 184.5
$ws.Range("L42").Value = 1783.2
$ws.Range("M42").Value = 45.5
$ws.Range("N42").Value = -2243.2

# Row 43
$ws.Range("H43").Value = 500
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 500
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").Value = 500
$ws.Range("N43").Value = -638

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0

# Row 46
$ws.Range("H46").Value = 922.5
$ws.Range("I46").Value = 650
$ws.Range("J46").Value = 1013.3333
$ws.Range("K46").Value = 1950
$ws.Range("L46").Value = 3039.9999
$ws.Range("M46").Value = -1831
$ws.Range("N46").Value = -3277.9999

# Row 48
$ws.Range("H48").Value = 1000
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 3000
$ws.Range("L48").ClearContents()
$ws.Range("N48").Value = 0
$ws.Range("M48").Value = -2708

# Row 51
$ws.Range("H51").Value = 1996.6666
$ws.Range("I51").Value = 1521.75
$ws.Range("J51").Value = 2169.3635
$ws.Range("K51").Value = 1521.75
$ws.Range("L51").Value = 2169.3635
$ws.Range("M51").Value = -1037.75
$ws.Range("N51").Value = -3137.3635

# Row 53
$ws.Range("H53").Value = 66897.664
$ws.Range("I53").Value = 500050.5
$ws.Range("J53").Value = 258.76923
$ws.Range("K53").Value = 500050.5
$ws.Range("L53").Value = 258.76923
$ws.Range("M53").Value = -499413.5
$ws.Range("N53").Value = -1532.76923

# Row 54
$ws.Range("H54").Value = 12800
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 12800
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 12800
$ws.Range("N54").Value = -13772

# Row 55
$ws.Range("H55").Value = 648.7692
$ws.Range("I55").Value = 1362
$ws.Range("J55").Value = 203
$ws.Range("K55").Value = 1362
$ws.Range("L55").Value = 203
$ws.Range("M55").Value = -1148
$ws.Range("N55").Value = -631

# Row 56
$ws.Range("H56").Value = 1000
$ws.Range("I56").Value = 1000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 3000
$ws.Range("L56").ClearContents()
$ws.Range("N56").Value = 0
$ws.Range("M56").Value = -2466

# Row 59
$ws.Range("H59").Value = 3999.6667
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 3999.6667
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 11999.0001
$ws.Range("N59").Value = -13113.0001

# Row 60
$ws.Range("H60").Value = 922.5
$ws.Range("I60").Value = 650
$ws.Range("J60").Value = 1013.3333
$ws.Range("K60").Value = 1950
$ws.Range("L60").Value = 3039.9999
$ws.Range("M60").Value = -1466
$ws.Range("N60").Value = -4007.9999

# Row 61
$ws.Range("H61").Value = 933.3333
$ws.Range("I61").Value = 150
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 450
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -278
$ws.Range("N61").Value = -7844

# Row 113
$ws.Range("H113").Value = 2392
$ws.Range("I113").Value = 2054
$ws.Range("J113").Value = 2899
$ws.Range("K113").Value = 2054
$ws.Range("L113").Value = 2899
$ws.Range("M113").Value = 1200
$ws.Range("N113").Value = -9407

# Row 133
$ws.Range("H133").Value = 50725
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50725
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50725
$ws.Range("N133").Value = -60845

# Row 134
$ws.Range("H134").Value = 58124.332
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 58124.332
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 58124.332
$ws.Range("N134").Value = -68264.33199999999

# Row 136
$ws.Range("H136").Value = 48990
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 48990
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 48990
$ws.Range("N136").Value = -59190

# Row 139
$ws.Range("H139").Value = 49126.668
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 49126.668
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 49126.668
$ws.Range("N139").Value = -59406.668

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 8773290
$ws.Range("I74").Value = 1095.7838
$ws.Range("J74").Value = 25001848
$ws.Range("K74").Value = 1095.7838
$ws.Range("L74").Value = 25001848
$ws.Range("M74").Value = -221.7837999999999
$ws.Range("N74").Value = -25003596

# Row 77
$ws.Range("H77").Value = 8773290
$ws.Range("I77").Value = 1095.7838
$ws.Range("J77").Value = 25001848
$ws.Range("K77").Value = 5478.919
$ws.Range("L77").Value = 125009240
$ws.Range("M77").Value = -1110.919
$ws.Range("N77").Value = -125017976

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1062.6216
$ws.Range("I94").Value = 783.7308
$ws.Range("J94").Value = 1721.8182
$ws.Range("K94").Value = 783.7308
$ws.Range("L94").Value = 1721.8182
$ws.Range("M94").Value = -332.7308
$ws.Range("N94").Value = -2623.8182

# Row 132
$ws.Range("H132").Value = 50262.855
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 50262.855
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 50262.855
$ws.Range("N132").Value = -60382.855

$ws = $wb.Worksheets.Item("CRP")
# Row 138
$ws.Range("H138").Value = 57637.777
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 57637.777
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 57637.777
$ws.Range("N138").Value = -67917.777

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1887966.1
$ws.Range("I131").Value = 8333774
$ws.Range("J131").Value = 1388.2927
$ws.Range("K131").Value = 25001322
$ws.Range("L131").Value = 4164.8781
$ws.Range("M131").Value = -24996282
$ws.Range("N131").Value = -14244.8781

$ws = $wb.Worksheets.Item("GSM")
# Row 135
$ws.Range("H135").Value = 55517.668
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 55517.668
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 55517.668
$ws.Range("N135").Value = -65657.66800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 16296714
$ws.Range("I122").Value = 35739290
$ws.Range("J122").Value = 3335000
$ws.Range("K122").Value = 107217870
$ws.Range("L122").Value = 10005000
$ws.Range("M122").Value = -107215420
$ws.Range("N122").Value = -10009900

# Row 132
$ws.Range("H132").Value = 23818282
$ws.Range("I132").Value = 27787218
$ws.Range("J132").Value = 4668
$ws.Range("K132").Value = 83361654
$ws.Range("L132").Value = 14004
$ws.Range("M132").Value = -83359124
$ws.Range("N132").Value = -19064

$ws = $wb.Worksheets.Item("WVR")
# Row 11
$ws.Range("H11").Value = 2000000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 2000000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 2000000
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -2000284

# Row 126
$ws.Range("H126").Value = 978.5
$ws.Range("I126").Value = 967.75
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 2903.25
$ws.Range("L126").Value = 1000
$ws.Range("M126").Value = -433.25
$ws.Range("N126").Value = -7940
